# Update "浏览/关注" (F column) counts that were scraped again at a later time.
# Values taken from commit "Update gh-pages to output generated at 456a3b4".

$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value  = 84
$ws.Range("F4").Value  = 1771
$ws.Range("F7").Value  = 1148
$ws.Range("F8").Value  = 1581
$ws.Range("F11").Value = 20
$ws.Range("F12").Value = 1506
$ws.Range("F13").Value = 3130
$ws.Range("F14").Value = 675
$ws.Range("F15").Value = 1823
$ws.Range("F16").Value = 1822
$ws.Range("F17").Value = 883
$ws.Range("F18").Value = 298
$ws.Range("F20").Value = 1508
$ws.Range("F21").Value = 305
$ws.Range("F23").Value = 24
$ws.Range("F24").Value = 1288
$ws.Range("F25").Value = 420
$ws.Range("F26").Value = 496
$ws.Range("F27").Value = 182
$ws.Range("F28").Value = 6773
$ws.Range("F29").Value = 5374
$ws.Range("F30").Value = 771
$ws.Range("F31").Value = 593
$ws.Range("F32").Value = 1707
$ws.Range("F34").Value = 228

# Sheet "演出"
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 89
$ws.Range("F7").Value = 106

# Sheet "全部类型"
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value  = 84
$ws.Range("F5").Value  = 89
$ws.Range("F7").Value  = 1771
$ws.Range("F10").Value = 1149
$ws.Range("F11").Value = 1581
$ws.Range("F15").Value = 20
$ws.Range("F16").Value = 1506
$ws.Range("F17").Value = 3130
$ws.Range("F18").Value = 675
$ws.Range("F19").Value = 1823
$ws.Range("F20").Value = 1822
$ws.Range("F21").Value = 883
$ws.Range("F22").Value = 298
$ws.Range("F24").Value = 1508
$ws.Range("F25").Value = 305
$ws.Range("F28").Value = 24
$ws.Range("F30").Value = 1288
$ws.Range("F31").Value = 420
$ws.Range("F32").Value = 496
$ws.Range("F33").Value = 182
$ws.Range("F34").Value = 6773
$ws.Range("F35").Value = 5374
$ws.Range("F36").Value = 771
$ws.Range("F37").Value = 593
$ws.Range("F38").Value = 1707
$ws.Range("F39").Value = 106
$ws.Range("F42").Value = 228
